# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the 0f3411bb-... row on the
# localization status report, reflecting a fresh handoff for both the
# zh-cn and de-de locale sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-29 03:37:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-29 03:38:09"
